$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-12 from 45185 to 45204
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
